$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "68.608.82"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "2.423.76"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue $ws "D5" "558.30"
$ws.Range("E5").Value = "  -0.39%  "
Set-TextValue $ws "D6" "161.10"
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("E7").Value = "  +0.00%  "
Set-TextValue $ws "D8" "0.513"
$ws.Range("E8").Value = "  +1.09%  "
Set-TextValue $ws "D9" "0.167"
$ws.Range("E9").Value = "  +10.13%  "
Set-TextValue $ws "D10" "0.162"
$ws.Range("E10").Value = "  -1.70%  "
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("E12").Value = "  -5.71%  "
$ws.Range("D13").Value = "68.525.68"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("E14").Value = "  +3.59%  "
$ws.Range("D15").Value = "2.872.50"
$ws.Range("E15").Value = "  -1.14%  "
Set-TextValue $ws "D16" "23.06"
$ws.Range("E16").Value = "  -2.46%  "
$ws.Range("D17").Value = "2.430.18"
$ws.Range("E17").Value = "  -1.94%  "
Set-TextValue $ws "D18" "10.44"
$ws.Range("E18").Value = "  -2.62%  "
Set-TextValue $ws "D19" "336.39"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("E22").Value = "  +1.58%  "
$ws.Range("E23").Value = "  -0.03%  "
Set-TextValue $ws "D24" "66.77"
$ws.Range("E24").Value = "  +0.03%  "
Set-TextValue $ws "D25" "3.67"
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("D26").Value = "2.556.30"
$ws.Range("E26").Value = "  -1.67%  "
Set-TextValue $ws "D27" "1.01"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("D29").Value = "0.0₃0812"
$ws.Range("E29").Value = "  -0.82%  "
Set-TextValue $ws "D30" "7.10"
$ws.Range("E30").Value = "  -1.69%  "
Set-TextValue $ws "D31" "1.00"
$ws.Range("E31").Value = "  +0.01%  "
Set-TextValue $ws "D32" "425.75"
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("E34").Value = "  -0.68%  "
Set-TextValue $ws "D35" "159.64"
$ws.Range("E35").Value = "  +0.39%  "
Set-TextValue $ws "D36" "19.01"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("E39").Value = "  -3.93%  "
Set-TextValue $ws "D40" "0.296"
$ws.Range("E40").Value = "  -1.49%  "
Set-TextValue $ws "D41" "1.49"
$ws.Range("E41").Value = "  +1.54%  "
Set-TextValue $ws "D42" "4.32"
$ws.Range("E42").Value = "  -2.64%  "
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("E44").Value = "  -2.55%  "
Set-TextValue $ws "D46" "130.74"
$ws.Range("E46").Value = "  -0.27%  "
Set-TextValue $ws "D47" "0.0715"
$ws.Range("E47").Value = "  +0.28%  "
Set-TextValue $ws "D48" "0.478"
$ws.Range("E48").Value = "  -1.30%  "
Set-TextValue $ws "D49" "0.554"
$ws.Range("E49").Value = "  -1.54%  "
Set-TextValue $ws "D50" "0.0920"
$ws.Range("E50").Value = "  +0.78%  "
Set-TextValue $ws "D51" "1.13"
